$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.309.61"
$ws.Range("E2").Value = "  -1.91%  "
$ws.Range("D3").Value = "1.835.33"
$ws.Range("E3").Value = "  -2.24%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "258.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -7.32%  "
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5193"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.26%  "
$ws.Range("E8").Value = "  -6.75%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06738"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.09%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.41"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -8.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7589"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.61%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07646"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.80%  "
$ws.Range("D13").Value = "1.816.62"
$ws.Range("E13").Value = "  -3.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.22"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.006"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.92%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.003"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.14%  "
$ws.Range("E17").Value = "  -4.49%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.002"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007864"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.61%  "
$ws.Range("D20").Value = "26.350.82"
$ws.Range("E20").Value = "  -1.92%  "
$ws.Range("D21").Value = "2.072.22"
$ws.Range("E21").Value = "  -2.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.537"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.48%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.389"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.29%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.908"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.57%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.06"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.230"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.92%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.642"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.26%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.90"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.81%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "111.23"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.155"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.66%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.130"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08698"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04763"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.86%  "
$ws.Range("E34").Value = "  -1.59%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.110"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.53%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6893"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.054"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01751"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.190"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -8.49%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4810"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "110.82"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.78%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.091"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8809"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.40%  "
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.606"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.81%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4111"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -8.62%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05845"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.64%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.974"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.85%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1227"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -8.74%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.55"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8806"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.49%  "
